# feat: add 2022-Q3 data
#
# Adds a new "2022-Q3" worksheet (created from a copy of the existing
# "2022-Q2" sheet, then repopulated with the new quarter's fund data),
# and records the new quarter's summary figures on the "总计" sheet.

function Set-TextValue($cell, $val) {
    # Forces the cell to hold a text value (matching the source data,
    # which stores numeric-looking figures like "1.20" as text), then
    # restores the cell to the default "Normal" style so no stray
    # number-format style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ2 = $wb.Worksheets.Item(2)

# --- 总计 (summary) sheet -------------------------------------------------
# Push the existing 2022-Q2 summary row down to row 3 (copy keeps the
# original formatting), then write the new 2022-Q3 summary into row 2.
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Cells.Item(3, 1).Value = 1

$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 0.08

# --- New 2022-Q3 sheet -----------------------------------------------------
# Duplicate the "2022-Q2" sheet (placed right after it) so the new sheet
# starts with the same layout/formatting, then swap the names so the
# original sheet (keeping its sheetId) becomes "2022-Q3" and the copy
# keeps the "2022-Q2" name with the old data untouched.
$wsQ2.Copy($null, $wsQ2)
$wsOldQ2Copy = $wb.Worksheets.Item(3)

$wsQ2.Name = "2022-Q3"
$wsOldQ2Copy.Name = "2022-Q2"

$wsQ3 = $wsQ2

Set-TextValue $wsQ3.Cells.Item(2, 4) "1.20"
Set-TextValue $wsQ3.Cells.Item(2, 5) "94.98"
Set-TextValue $wsQ3.Cells.Item(2, 6) "4.34"
Set-TextValue $wsQ3.Cells.Item(2, 7) "0.0521"

Set-TextValue $wsQ3.Cells.Item(3, 4) "0.60"
Set-TextValue $wsQ3.Cells.Item(3, 5) "94.98"
Set-TextValue $wsQ3.Cells.Item(3, 6) "4.34"
Set-TextValue $wsQ3.Cells.Item(3, 7) "0.0260"

# Match the "总计" sheet's look (header/row formatting and page margins)
# for the brand-new 2022-Q3 sheet.
$wsTotal.Range("B1:D1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A3").PasteSpecial(-4122)

$wsQ3.PageSetup.LeftMargin = 54
$wsQ3.PageSetup.RightMargin = 54
$wsQ3.PageSetup.TopMargin = 72
$wsQ3.PageSetup.BottomMargin = 72
$wsQ3.PageSetup.HeaderMargin = 36
$wsQ3.PageSetup.FooterMargin = 36

Write-Host "Added 2022-Q3 sheet and updated 总计 summary."
